$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 5373.000859812833
$ws.Range("D4").Value = 3268.002014659345
$ws.Range("E4").Value = 0.6082265943976016
$ws.Range("F4").Value = 1.644124096530862
$ws.Range("G4").Value = 145.9361429969722
$ws.Range("H4").Value = 15.89677538211254
$ws.Range("I4").Value = 12.25906828093957
$ws.Range("J4").Value = 290.2080667279661
$ws.Range("K4").Value = 280.3078072129283
$ws.Range("L4").Value = 163.4015853721648
$ws.Range("M4").Value = 2.784134542991524
$ws.Range("N4").Value = 4.935462251400168
$ws.Range("O4").Value = 0.007548424894594064
$ws.Range("P4").Value = 348.7826346189249
$ws.Range("Q4").Value = 522.286516635213
$ws.Range("R4").Value = 21.98354299389757
$ws.Range("S4").Value = -2915.938725396452
$ws.Range("T4").Value = -313.0000453908506
$ws.Range("U4").Value = -245.1738171938969
$ws.Range("V4").Value = -5455.378699940396
$ws.Range("W4").Value = -2745.745190808084
$ws.Range("X4").Value = -5584.172601264669
